# SRS updated to version 0.4 - version with Mark's input and all the classes added.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: "Update Queries" assignment moved from "all" to "Carlos", status set to "Finished"
$ws.Range("F4").Value = "Carlos"
$ws.Range("G4").Value = "Finished"

# Row 5: "Exception handling" assignment moved from "all" to "Sanjay", status set to "Finished"
$ws.Range("F5").Value = "Sanjay"
$ws.Range("G5").Value = "Finished"

# Row 15: "UI" assignment moved from "all" to "Sanjay - Carlos", status set to "Finished"
$ws.Range("F15").Value = "Sanjay - Carlos"
$ws.Range("G15").Value = "Finished"

# Row 16: "Test" status set to "Finished"
$ws.Range("G16").Value = "Finished"

# Update the current selection to match the author's last cursor position
[void]$ws.Range("I13").Select()
